# Add a new account row (VINICIUS, 005886225, 5000) to the "Export" sheet,
# inserted directly above the existing "004313254" (GUSTAVO) row — i.e. as
# the new row 11 (header is row 1; CARLA/004643153 stays row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift GUSTAVO (and everything below) down by one row.
$ws.Rows.Item(11).Insert()

# Account numbers are zero-padded strings, so force text formatting before
# assigning the value — otherwise Excel would parse "005886225" as the
# number 5886225 and drop the leading zeros.
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "005886225"
$ws.Cells.Item(11, 2).Value = "VINICIUS"
$ws.Cells.Item(11, 3).Value = 5000
